# Append a new scrape run (2025-09-20 18:22:20 JST) to the "ランサーズ" sheet.
# - bumps every existing row's timestamp to the new scrape time
# - inserts 3 brand-new job rows ahead of the previously-first "new" job
# - widens column B to fit the longer titles
# - rebuilds the hyperlinks collection so F2:F12 all point at the right URL

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-09-20 18:22:20"

# --- 1. Make room: push the old rows 7-9 down to rows 10-12 ----------------
$ws.Rows("7:9").Insert()

# --- 2. Refresh the timestamp on the rows that were already present --------
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 3. Two brand-new postings land in the freshly inserted rows 7 & 8 -----
$ws.Cells.Item(7, 1).Value = $newTimestamp
$ws.Cells.Item(7, 2).Value = "【急募】WordPress開発者を探しています!魅力的な案件です"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5397452"
$ws.Cells.Item(7, 7).Value = 88
$ws.Cells.Item(7, 8).Value = "◆開発 ○WordPress"

$ws.Cells.Item(8, 1).Value = $newTimestamp
$ws.Cells.Item(8, 2).Value = "【急募】LARAVEL 開発者を募集しています!高報酬案件"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5397425"
$ws.Cells.Item(8, 7).Value = 68
$ws.Cells.Item(8, 8).Value = "◆開発"

# --- 4. Row 9 keeps the old "教育系" posting, just with the new timestamp --
$ws.Cells.Item(9, 1).Value = $newTimestamp
$ws.Cells.Item(9, 2).Value = "【急募】教育系のWEBサイトの作成"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5397264"
$ws.Cells.Item(9, 7).Value = 33
$ws.Cells.Item(9, 8).Value = "◇サイト"

# --- 5. A third brand-new posting in the new row 10 -------------------------
$ws.Cells.Item(10, 1).Value = $newTimestamp
$ws.Cells.Item(10, 2).Value = "GoogleスプレッドシートへMYSQLデータ取り込み及びスプレッドシート改修"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5397403"
$ws.Cells.Item(10, 7).Value = 30
$ws.Cells.Item(10, 8).Value = "◇MySQL"

# --- 6. Row 11 keeps the old "フォートナイト" posting, new timestamp -------
$ws.Cells.Item(11, 1).Value = $newTimestamp
$ws.Cells.Item(11, 2).Value = "【フォートナイト】クリエイティブ作品を世界に公開したい!"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5397192"
$ws.Cells.Item(11, 7).Value = 18

# --- 7. Row 12 keeps the old "Web広告" posting, new timestamp --------------
$ws.Cells.Item(12, 1).Value = $newTimestamp
$ws.Cells.Item(12, 2).Value = "初回 Web広告のタグ設置・動作確認"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5397007"
$ws.Cells.Item(12, 7).Value = 18

# --- 8. Widen column B so the longer titles fit -----------------------------
$ws.Columns("B").ColumnWidth = 40.17

# --- 9. Rebuild the hyperlinks top to bottom (row-insert doesn't move them) -
$ws.Cells.Hyperlinks.Delete()
for ($r = 2; $r -le 12; $r++) {
    $target = $ws.Cells.Item($r, 6).Value2
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $target)
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
